$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shift existing Week #13 rows down to make room for the two new
#     Week #14 entries (Bettis, Dahl) that are inserted at the top of the
#     injury list, then refresh the two still-active players' notes/dates,
#     and append Bryan Shaw as a brand new entry. ---

# Row 2: Chad Bettis (new entry)
$ws.Range("A2").Value = "Chad Bettis"
$ws.Range("B2").Value = "bettich01"
$ws.Range("C2").Value = "July 01 2018"
$ws.Range("D2").Value = "Finger"
$ws.Range("E2").Value = "Bettis left his last outing due to a right middle finger injury and it is unknown if he will make his next scheduled start Saturday against the Mariners."

# Row 3: David Dahl (new entry)
$ws.Range("A3").Value = "David Dahl"
$ws.Range("B3").Value = "dahlda01"
$ws.Range("C3").Value = "June 02 2018"
$ws.Range("D3").Value = "Foot"
$ws.Range("E3").Value = "Dahl has been placed on the 10-day disabled list with a broken right foot and will likely need six-to-eight weeks to recover."

# Row 4: Carlos Estevez (existing player, moved down, note/date refreshed)
$ws.Range("A4").Value = "Carlos Estevez"
$ws.Range("B4").Value = "estevca01"
$ws.Range("C4").Value = "July 01 2018"
$ws.Range("D4").Value = "Oblique"
$ws.Range("E4").Value = "Estevez has been shifted to the 60-day disabled list due to a left oblique strain but is likely to return during the middle portion of July."
$ws.Rows.Item(4).RowHeight = 30
$ws.Range("E4").WrapText = $true

# Row 5: Zac Rosscup (existing player, moved down, note/date refreshed)
$ws.Range("A5").Value = "Zac Rosscup"
$ws.Range("B5").Value = "rosscza01"
$ws.Range("C5").Value = "July 01 2018"
$ws.Range("D5").Value = "Finger"
$ws.Range("E5").Value = "Rosscup has been shifted to the 60-day disabled list dealing with warts on his left middle finger and is likely to stay on the DL until at least the middle of July."
$ws.Rows.Item(5).RowHeight = 30
$ws.Range("E5").WrapText = $true

# Row 6: Bryan Shaw (new entry)
$ws.Range("A6").Value = "Bryan Shaw"
$ws.Range("B6").Value = "shawbr01"
$ws.Range("C6").Value = "June 24 2018"
$ws.Range("D6").Value = "Calf"
$ws.Range("E6").Value = "Shaw has landed on the 10-day disables list with a right calf strain and there is no timetable for return."
$ws.Rows.Item(6).RowHeight = 30
$ws.Range("E6").WrapText = $true

# Update the selection left behind on the sheet after the edits
$ws.Range("E16").Select()
